$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.772.81"
$ws.Range("E2").Value = "  -2.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.782.40"
$ws.Range("E3").Value = "  -2.09%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'310.62"
$ws.Range("E5").Value = "  -2.11%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.06%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.5109"
$ws.Range("E7").Value = "  -1.10%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3823"
$ws.Range("E8").Value = "  -1.66%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07795"
$ws.Range("E9").Value = "  -7.93%  "

# Row 10 - was OKB, now Polygon (rows 10 and 11 swapped)
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.086"
$ws.Range("E10").Value = "  -2.31%  "

# Row 11 - was Polygon, now OKB
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'40.67"
$ws.Range("E11").Value = "  -2.77%  "

# Row 12 - BinanceUSD
$ws.Range("E12").Value = "  -0.04%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'6.189"
$ws.Range("E13").Value = "  -3.98%  "

# Row 14 - Solana
$ws.Range("E14").Value = "  -4.51%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "1.781.33"
$ws.Range("E15").Value = "  -1.90%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'7.186"

# Row 17 - Litecoin
$ws.Range("D17").Value = "'91.34"
$ws.Range("E17").Value = "  -1.61%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -5.99%  "

# Row 19 - TRON
$ws.Range("D19").Value = "'0.06548"
$ws.Range("E19").Value = "  -1.28%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.02%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "'16.99"
$ws.Range("E21").Value = "  -4.31%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.908"
$ws.Range("E22").Value = "  -3.03%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "27.822.71"
$ws.Range("E23").Value = "  -2.46%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  -4.14%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.236"
$ws.Range("E25").Value = "  -1.71%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'159.94"
$ws.Range("E26").Value = "  +0.42%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -4.18%  "

# Row 28 - WrappedliquidstakedEther2.0
$ws.Range("D28").Value = "1.986.54"
$ws.Range("E28").Value = "  -1.93%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.358"
$ws.Range("E29").Value = "  -1.61%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'123.71"
$ws.Range("E30").Value = "  -1.44%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.1070"
$ws.Range("E31").Value = "  -1.82%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'1.032"
$ws.Range("E32").Value = "  -5.77%  "

# Row 33 - HuobiToken
$ws.Range("E33").Value = "  -0.33%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  -4.37%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "'0.07073"
$ws.Range("E35").Value = "  -4.96%  "

# Row 36 - FraxShare
$ws.Range("E36").Value = "  -0.31%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02303"
$ws.Range("E37").Value = "  -2.42%  "

# Row 38 - Algorand
$ws.Range("D38").Value = "'0.2118"
$ws.Range("E38").Value = "  -5.13%  "

# Row 39 - Aptos
$ws.Range("D39").Value = "'11.47"
$ws.Range("E39").Value = "  +1.74%  "

# Row 40 - InternetComputer(DFINITY)
$ws.Range("D40").Value = "'4.987"
$ws.Range("E40").Value = "  -4.28%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "'0.6083"
$ws.Range("E41").Value = "  -3.66%  "

# Row 42 - Frax
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  +0.01%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "'1.153"
$ws.Range("E43").Value = "  -3.19%  "

# Row 44 - WEMIXTOKEN
$ws.Range("D44").Value = "'1.323"
$ws.Range("E44").Value = "  -5.45%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "'13.11"
$ws.Range("E45").Value = "  -3.56%  "

# Row 46 - was Decentraland, now PancakeSwap (rows 46 and 47 swapped)
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.707"
$ws.Range("E46").Value = "  -2.08%  "

# Row 47 - was PancakeSwap, now Decentraland
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5872"
$ws.Range("E47").Value = "  -1.26%  "

# Row 48 - Quant
$ws.Range("D48").Value = "'125.81"
$ws.Range("E48").Value = "  -0.48%  "

# Row 49 - EOS
$ws.Range("D49").Value = "'1.195"
$ws.Range("E49").Value = "  -0.66%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "'1.895"
$ws.Range("E50").Value = "  -4.68%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "'0.06865"
$ws.Range("E51").Value = "  -1.66%  "
